$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Anabaena variabilis PCC 7120 DSM 107007"
$ws.Range("A5").Value = "Bacillus cereus NRS 248 ATCC 10987"
$ws.Range("A6").Value = "ZymoGut"
$ws.Range("A7").Value = "ZymoOral"
$ws.Range("A8").Value = "ZymoFecal"

$ws.Range("A9").Select()
